$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark: it currently sits right after the inline
#    picture ("截图20190703000812"); it should instead sit in the middle of
#    the paragraph that talks about web.config <modules>, right after
#    "...下面我" (i.e. right before "们就在WEB项目工程文档下的...").
#    Re-adding a bookmark named "_GoBack" moves the (unique) bookmark, which
#    also naturally splits the run at that text position - matching how
#    Word leaves the cursor mark after the last edit.
$splitPoint = $d.Content
$splitPoint.Find.Execute("下面我", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($splitPoint.End, $splitPoint.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# 2) Update the HttpHandler paragraph to also mention MvcHandler.
$old1 = "像一般处理程序 *.ashx、WebForm Page它们都是属于"
$new1 = "像一般处理程序 *.ashx、WebForm Page和MvcHandler它们都是属于"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "一般处理程序和Page其实都是继承"
$new2 = "一般处理程序、Page和MvcHandler其实都是继承"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Add the built-in "Strong" character style (based on Default Paragraph
#    Font, bold run formatting) to the style sheet, the same way Word mints
#    it into styles.xml the first time it is referenced - without leaving
#    any visible trace on the body content, apply/unapply it on a scratch
#    paragraph that is inserted and removed again.
$scratch = $d.Content
$scratch.Collapse(0)
$scratch.InsertParagraphAfter()
$scratchPara = $d.Paragraphs($d.Paragraphs.Count)
$scratchPara.Range.Text = "x"
$scratchRun = $d.Range($scratchPara.Range.Start, $scratchPara.Range.Start + 1)
$scratchRun.Style = "Strong"
$scratchPara.Range.Delete()

$strongStyle = $d.Styles("Strong")
$strongStyle.Priority = 0
